# Auto-generated Excel COM-interop edit script
# Applies numeric updates to leve-profit calculation columns (H-N)
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1156.2222
$ws.Range("I40").Value = 1175.125
$ws.Range("K40").Value = 1175.125
$ws.Range("M40").Value = -1000.125
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = $null
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").Value = $null
$ws.Range("H64").Value = 447348.8
$ws.Range("I64").Value = 853636.7
$ws.Range("J64").Value = 4125.727
$ws.Range("K64").Value = 853636.7
$ws.Range("L64").Value = 4125.727
$ws.Range("M64").Value = -853388.7
$ws.Range("N64").Value = -4621.727
$ws.Range("H67").Value = 447348.8
$ws.Range("I67").Value = 853636.7
$ws.Range("J67").Value = 4125.727
$ws.Range("K67").Value = 853636.7
$ws.Range("L67").Value = 4125.727
$ws.Range("M67").Value = -852778.7
$ws.Range("N67").Value = -5841.727
$ws.Range("H107").Value = 933.6667
$ws.Range("I107").Value = 945.9091
$ws.Range("J107").Value = 799
$ws.Range("K107").Value = 945.9091
$ws.Range("L107").Value = 799
$ws.Range("M107").Value = 974.0909
$ws.Range("N107").Value = -4639
$ws.Range("H125").Value = 5727.1
$ws.Range("I125").Value = 322.33334
$ws.Range("J125").Value = 6680.8823
$ws.Range("K125").Value = 2901.00006
$ws.Range("L125").Value = 60127.9407
$ws.Range("M125").Value = -441.0000600000003
$ws.Range("N125").Value = -65047.9407
$ws.Range("H135").Value = 68183060
$ws.Range("I135").Value = 35715860
$ws.Range("J135").Value = 125000650
$ws.Range("K135").Value = 321442740
$ws.Range("L135").Value = 1125005850
$ws.Range("M135").Value = -321440205
$ws.Range("N135").Value = -1125010920
$ws.Range("H137").Value = 2119.681
$ws.Range("I137").Value = 1846.7179
$ws.Range("K137").Value = 5540.153700000001
$ws.Range("M137").Value = -2990.153700000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6624.9443
$ws.Range("I61").Value = 6082.737
$ws.Range("J61").Value = 7230.9414
$ws.Range("K61").Value = 6082.737
$ws.Range("L61").Value = 7230.9414
$ws.Range("M61").Value = -5870.737
$ws.Range("N61").Value = -7654.9414
$ws.Range("H74").Value = 5744.5107
$ws.Range("I74").Value = 4265.421
$ws.Range("J74").Value = 11989.556
$ws.Range("K74").Value = 4265.421
$ws.Range("L74").Value = 11989.556
$ws.Range("M74").Value = -3391.421
$ws.Range("N74").Value = -13737.556
$ws.Range("H76").Value = 38818.285
$ws.Range("J76").Value = 38818.285
$ws.Range("L76").Value = 38818.285
$ws.Range("N76").Value = -39494.285
$ws.Range("H77").Value = 5744.5107
$ws.Range("I77").Value = 4265.421
$ws.Range("J77").Value = 11989.556
$ws.Range("K77").Value = 21327.105
$ws.Range("L77").Value = 59947.78
$ws.Range("M77").Value = -16959.105
$ws.Range("N77").Value = -68683.78
$ws.Range("H79").Value = 38818.285
$ws.Range("J79").Value = 38818.285
$ws.Range("L79").Value = 38818.285
$ws.Range("N79").Value = -41158.285
$ws.Range("H96").Value = 50000
$ws.Range("J96").Value = 50000
$ws.Range("L96").Value = 50000
$ws.Range("N96").Value = -55492
$ws.Range("H102").Value = 3048.3333
$ws.Range("I102").Value = 2672.5
$ws.Range("J102").Value = 3800
$ws.Range("K102").Value = 2672.5
$ws.Range("L102").Value = 3800
$ws.Range("M102").Value = -1050.5
$ws.Range("N102").Value = -7044
$ws.Range("H136").Value = 6624.9443
$ws.Range("I136").Value = 6082.737
$ws.Range("J136").Value = 7230.9414
$ws.Range("K136").Value = 18248.211
$ws.Range("L136").Value = 21692.8242
$ws.Range("M136").Value = -15698.211
$ws.Range("N136").Value = -26792.8242

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2537.4
$ws.Range("I134").Value = 2451.1667
$ws.Range("J134").Value = 2882.3333
$ws.Range("K134").Value = 7353.500100000001
$ws.Range("L134").Value = 8646.999899999999
$ws.Range("M134").Value = -4818.500100000001
$ws.Range("N134").Value = -13716.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4291.7393
$ws.Range("J31").Value = 3133.238
$ws.Range("L31").Value = 3133.238
$ws.Range("N31").Value = -3723.238
$ws.Range("H34").Value = 4291.7393
$ws.Range("J34").Value = 3133.238
$ws.Range("L34").Value = 3133.238
$ws.Range("N34").Value = -3537.238
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 854.63635
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 3000
$ws.Range("N4").Value = -3224
$ws.Range("H5").Value = 7936963.5
$ws.Range("I5").Value = 461.1579
$ws.Range("J5").Value = 83333736
$ws.Range("K5").Value = 1383.4737
$ws.Range("L5").Value = 250001208
$ws.Range("M5").Value = -1271.4737
$ws.Range("N5").Value = -250001432
$ws.Range("H34").Value = 2361.05
$ws.Range("J34").Value = 3077.7334
$ws.Range("L34").Value = 9233.200199999999
$ws.Range("N34").Value = -9401.200199999999
$ws.Range("H39").Value = 6545.7915
$ws.Range("I39").Value = 1400
$ws.Range("J39").Value = 7280.905
$ws.Range("K39").Value = 4200
$ws.Range("L39").Value = 21842.715
$ws.Range("M39").Value = -3906
$ws.Range("N39").Value = -22430.715
$ws.Range("H55").Value = 4250
$ws.Range("J55").Value = 4250
$ws.Range("L55").Value = 12750
$ws.Range("N55").Value = -13104
$ws.Range("H59").Value = 2000
$ws.Range("I59").Value = 1250
$ws.Range("J59").Value = 3500
$ws.Range("K59").Value = 3750
$ws.Range("L59").Value = 10500
$ws.Range("M59").Value = -3210
$ws.Range("N59").Value = -11580
$ws.Range("H113").Value = 682.42426
$ws.Range("J113").Value = 654.8276
$ws.Range("L113").Value = 1964.4828
$ws.Range("N113").Value = -6304.4828
$ws.Range("H131").Value = 31177.188
$ws.Range("I131").Value = 1419.0625
$ws.Range("J131").Value = 60935.312
$ws.Range("K131").Value = 4257.1875
$ws.Range("L131").Value = 182805.936
$ws.Range("M131").Value = 782.8125
$ws.Range("N131").Value = -192885.936
$ws.Range("H132").Value = 2025.1875
$ws.Range("I132").Value = 2145.4
$ws.Range("J132").Value = 1824.8334
$ws.Range("K132").Value = 19308.6
$ws.Range("L132").Value = 16423.5006
$ws.Range("M132").Value = -16778.6
$ws.Range("N132").Value = -21483.5006
$ws.Range("H135").Value = 7936963.5
$ws.Range("I135").Value = 461.1579
$ws.Range("J135").Value = 83333736
$ws.Range("K135").Value = 4150.4211
$ws.Range("L135").Value = 750003624
$ws.Range("M135").Value = -1615.4211
$ws.Range("N135").Value = -750008694

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 45.285713
$ws.Range("I2").Value = 48.333332
$ws.Range("J2").Value = 43
$ws.Range("K2").Value = 48.333332
$ws.Range("L2").Value = 43
$ws.Range("M2").Value = 64.666668
$ws.Range("N2").Value = -269
$ws.Range("H82").Value = 31507.857
$ws.Range("J82").Value = 31507.857
$ws.Range("L82").Value = 31507.857
$ws.Range("N82").Value = -32273.857
$ws.Range("H85").Value = 31507.857
$ws.Range("J85").Value = 31507.857
$ws.Range("L85").Value = 31507.857
$ws.Range("N85").Value = -34159.857
$ws.Range("H97").Value = 1590.1052
$ws.Range("I97").Value = 2248.3333
$ws.Range("J97").Value = 1286.3077
$ws.Range("K97").Value = 2248.3333
$ws.Range("L97").Value = 1286.3077
$ws.Range("M97").Value = -1752.3333
$ws.Range("N97").Value = -2278.3077
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1800
$ws.Range("I46").Value = 3000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2812
$ws.Range("H98").Value = 29355
$ws.Range("J98").Value = 29355
$ws.Range("L98").Value = 29355
$ws.Range("N98").Value = -35345
$ws.Range("H100").Value = 4140.4116
$ws.Range("I100").Value = 3534.7856
$ws.Range("J100").Value = 6966.6665
$ws.Range("K100").Value = 3534.7856
$ws.Range("L100").Value = 6966.6665
$ws.Range("M100").Value = -2993.7856
$ws.Range("N100").Value = -8048.6665
$ws.Range("H122").Value = 6124.25
$ws.Range("I122").Value = 5320.2163
$ws.Range("J122").Value = 7690
$ws.Range("K122").Value = 15960.6489
$ws.Range("L122").Value = 23070
$ws.Range("M122").Value = -13510.6489
$ws.Range("N122").Value = -27970
$ws.Range("H132").Value = 4211.115
$ws.Range("I132").Value = 4649.5713
$ws.Range("J132").Value = 3308.4119
$ws.Range("K132").Value = 13948.7139
$ws.Range("L132").Value = 9925.235700000001
$ws.Range("M132").Value = -11418.7139
$ws.Range("N132").Value = -14985.2357
$ws.Range("H136").Value = 4539.2554
$ws.Range("I136").Value = 2510.577
$ws.Range("J136").Value = 7050.952
$ws.Range("K136").Value = 7531.731000000001
$ws.Range("L136").Value = 21152.856
$ws.Range("M136").Value = -4981.731000000001
$ws.Range("N136").Value = -26252.856
$ws.Range("H140").Value = 56808
$ws.Range("J140").Value = 56808
$ws.Range("L140").Value = 56808
$ws.Range("N140").Value = -67168

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 25000
$ws.Range("I92").Value = 25000
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 25000
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -22504
$ws.Range("N92").Value = $null
$ws.Range("H132").Value = 1567.9298
$ws.Range("I132").Value = 815.37836
$ws.Range("J132").Value = 2960.15
$ws.Range("K132").Value = 2446.13508
$ws.Range("L132").Value = 8880.450000000001
$ws.Range("M132").Value = 83.86491999999998
$ws.Range("N132").Value = -13940.45
